$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns D (Price) and E (Volume(1h)) hold numeric-looking / percent-looking
# *text* values (t="inlineStr" in the source, General-formatted, no leading
# apostrophe). Force text storage via NumberFormat "@" before assigning so
# Excel does not auto-coerce "295.88" / "3.03%" into real numbers/percentages,
# then restore "General" display formatting once the literal text is stored.
$textCells = @(
  'D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51'
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# --- Apply the scraped values (Coin / Link / Price / Volume(1h)) ---
$ws.Range('D2').Value = '295.88'
$ws.Range('E2').Value = '3.03%'
$ws.Range('D3').Value = '41.32'
$ws.Range('E3').Value = '3.04%'
$ws.Range('D4').Value = '5.041'
$ws.Range('E4').Value = '0.07%'
$ws.Range('D5').Value = '0.07461'
$ws.Range('E5').Value = '2.06%'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Value = '1.574'
$ws.Range('E6').Value = '1.71%'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Value = '0.9341'
$ws.Range('E7').Value = '2.72%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '2.402'
$ws.Range('E8').Value = '0.21%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = '0.1198'
$ws.Range('E9').Value = '-0.13%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1812'
$ws.Range('E10').Value = '4.16%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.08809'
$ws.Range('E11').Value = '1.61%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.04349'
$ws.Range('E12').Value = '4.55%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.1050'
$ws.Range('E13').Value = '-0.13%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001263'
$ws.Range('E14').Value = '-0.88%'
$ws.Range('D15').Value = '0.005869'
$ws.Range('E15').Value = '1.07%'
$ws.Range('E16').Value = '-1.19%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '4.353'
$ws.Range('E17').Value = '1.63%'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').Value = '0.3306'
$ws.Range('E18').Value = '0.70%'
$ws.Range('B19').Value = 'MCDex'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D19').Value = '8.023'
$ws.Range('E19').Value = '6.09%'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').Value = '0.1378'
$ws.Range('E20').Value = '2.64%'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D21').Value = '0.2972'
$ws.Range('E21').Value = '2.96%'
$ws.Range('D22').Value = '0.04015'
$ws.Range('E22').Value = '4.55%'
$ws.Range('D23').Value = '0.001268'
$ws.Range('E23').Value = '-0.20%'
$ws.Range('D24').Value = '0.003874'
$ws.Range('E24').Value = '5.04%'
$ws.Range('D25').Value = '0.0001233'
$ws.Range('E25').Value = '-4.01%'
$ws.Range('D26').Value = '0.0003729'
$ws.Range('E26').Value = '-0.15%'
$ws.Range('D38').Value = '0.02376'
$ws.Range('E38').Value = '2.29%'
$ws.Range('D39').Value = '0.05171'
$ws.Range('E39').Value = '3.61%'
$ws.Range('D40').Value = '0.006033'
$ws.Range('E40').Value = '18.08%'
$ws.Range('D41').Value = '0.007799'
$ws.Range('E41').Value = '1.23%'
$ws.Range('D42').Value = '0.1318'
$ws.Range('E42').Value = '3.81%'
$ws.Range('D43').Value = '0.007388'
$ws.Range('E43').Value = '-0.23%'
$ws.Range('D44').Value = '0.007853'
$ws.Range('E44').Value = '4.19%'
$ws.Range('D45').Value = '0.2939'
$ws.Range('E45').Value = '-5.37%'
$ws.Range('D46').Value = '0.00006355'
$ws.Range('E46').Value = '-2.65%'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').Value = '-0.14%'
$ws.Range('D48').Value = '0.04644'
$ws.Range('E48').Value = '-81.56%'
$ws.Range('D49').Value = '0.004207'
$ws.Range('E49').Value = '0.14%'
$ws.Range('D50').Value = '0.00002104'
$ws.Range('E50').Value = '-0.14%'
$ws.Range('D51').Value = '0.0002004'
$ws.Range('E51').Value = '-0.14%'

foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "General"
}
